$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# Force text entry so numeric-looking strings (e.g. "1.00", "8.40") keep their
# exact formatting instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.786.08"
$ws.Range("D3").Value = "2.992.55"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "542.18"
$ws.Range("D6").Value = "153.13"
$ws.Range("D9").Value = "3.000.31"
$ws.Range("D13").Value = "3.515.11"
$ws.Range("D14").Value = "0.123"
$ws.Range("D15").Value = "61.893.32"
$ws.Range("D17").Value = "2.993.04"
$ws.Range("D19").Value = "390.15"
$ws.Range("D21").Value = "11.91"
$ws.Range("D22").Value = "6.61"
$ws.Range("D24").Value = "64.93"
$ws.Range("D25").Value = "0.467"
$ws.Range("D27").Value = "0.998"
$ws.Range("D28").Value = "0.0₃0941"
$ws.Range("D29").Value = "8.40"
$ws.Range("D31").Value = "1.71"
$ws.Range("D32").Value = "20.40"
$ws.Range("D33").Value = "159.54"
$ws.Range("D35").Value = "4.61"
$ws.Range("D38").Value = "1.57"
$ws.Range("D39").Value = "2.432.88"
$ws.Range("D40").Value = "22.39"
$ws.Range("D42").Value = "37.18"
$ws.Range("D47").Value = "4.93"
$ws.Range("D48").Value = "0.0953"
$ws.Range("D49").Value = "19.65"
$ws.Range("D50").Value = "10.45"
$ws.Range("D51").Value = "264.85"

# Restore the default (unstyled) look for column D now that the text is set.
$ws.Range("D2:D51").Style = "Normal"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -5.68%  "
$ws.Range("E3").Value = "  -6.65%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -6.02%  "
$ws.Range("E6").Value = "  -8.82%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -5.57%  "
$ws.Range("E10").Value = "  -6.98%  "
$ws.Range("E11").Value = "  -7.54%  "
$ws.Range("E12").Value = "  -6.84%  "
$ws.Range("E13").Value = "  -6.55%  "
$ws.Range("E14").Value = "  -3.86%  "
$ws.Range("E15").Value = "  -5.41%  "
$ws.Range("E16").Value = "  -8.20%  "
$ws.Range("E17").Value = "  -7.20%  "
$ws.Range("E18").Value = "  -7.19%  "
$ws.Range("E19").Value = "  -5.81%  "
$ws.Range("E20").Value = "  -4.55%  "
$ws.Range("E21").Value = "  -7.49%  "
$ws.Range("E22").Value = "  -8.03%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("E24").Value = "  -6.66%  "
$ws.Range("E25").Value = "  -4.64%  "
$ws.Range("E26").Value = "  -8.12%  "
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("E28").Value = "  -11.22%  "
$ws.Range("E29").Value = "  -6.03%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -7.28%  "
$ws.Range("E32").Value = "  -5.88%  "
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("E34").Value = "  -6.39%  "
$ws.Range("E35").Value = "  -8.66%  "
$ws.Range("E36").Value = "  -6.95%  "
$ws.Range("E37").Value = "  -6.89%  "
$ws.Range("E38").Value = "  -10.14%  "
$ws.Range("E39").Value = "  -10.93%  "
$ws.Range("E40").Value = "  -7.63%  "
$ws.Range("E41").Value = "  -6.93%  "
$ws.Range("E42").Value = "  -4.94%  "
$ws.Range("E43").Value = "  -7.25%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("E46").Value = "  -7.09%  "
$ws.Range("E47").Value = "  -11.95%  "
$ws.Range("E48").Value = "  -3.89%  "
$ws.Range("E49").Value = "  -9.52%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("E51").Value = "  -10.97%  "

# --- Column B (Coin) / C (Link) swap for rows 40 and 42 ---
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("B42").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
